$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.51305890083313
$ws.Range("B1").Value = 1.160800695419312
$ws.Range("C1").Value = 1.575078129768372
$ws.Range("D1").Value = 2.389882326126099
$ws.Range("E1").Value = 6.440203666687012
